# GrowingShadeUseReport.docx content refresh ("make it work again")
#
# 1. Bump the report date by a day.
# 2. Figure 1.1 was mislabeled with the Figure 1.2 description/caption text
#    ("Average time that each user has spent engaged...") - fix it to
#    describe what the figure actually shows (users by city).
# 3. Clarify the "highest engagement" sentence in the application-usage
#    section.

$d = $word.ActiveDocument

# 1. Report date.
$d.Content.Find.Execute("17 March 2022", $true, $false, $false, $false, $false, `
    $true, 1, $false, "18 March 2022", 2) | Out-Null

# 2a. Figure 1.1 accessible description (wp:docPr/@descr) - only the first
#     inline picture in the document.
$d.InlineShapes.Item(1).AlternativeText = `
    "Figure 1.1: Users accessing the Growing Shade application by city."

# 2b. Figure 1.1 visible caption paragraph underneath the picture.
$d.Content.Find.Execute( `
    "Figure 1.1: Average time that each user has spent engaged on the Growing Shade application.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Figure 1.1: Users accessing the Growing Shade application by city.", 2) | Out-Null

# 3. "Highest engagement" sentence: add the "on average" qualifier and
#    reword "on" -> "within" the application.
$d.Content.Find.Execute( `
    "have spent the most time on the application", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "have spent the most time, on average, within the application", 2) | Out-Null
